$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Write-Host "I1:" $ws.Range("I1").Value()
Write-Host "I3:" $ws.Range("I3").Value()
Write-Host "I4:" $ws.Range("I4").Value()
Write-Host "I8:" $ws.Range("I8").Value()
Write-Host "D3:" $ws.Range("D3").Value()
Write-Host "D2:" $ws.Range("D2").Value()
Write-Host "D4:" $ws.Range("D4").Value()
